$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "29.270.37"
$ws.Range("E2").Value = "  -0.42%  "

$ws.Range("D3").Value = "1.830.78"
$ws.Range("E3").Value = "  -0.61%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.54%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.79%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6025"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.99%  "

$ws.Range("E7").Value = "  +0.46%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07025"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.66%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2793"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.56%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.66%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07661"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.62%  "

$ws.Range("D12").Value = "1.833.87"
$ws.Range("E12").Value = "  -0.55%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.799"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.000009886"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.81%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6251"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "79.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.24%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "29.275.87"
$ws.Range("E17").Value = "  -0.51%  "

$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.839"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.51%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "224.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.63%  "

$ws.Range("E20").Value = "  +0.44%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.97%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.005"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.006"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "156.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.998"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.85%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1299"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.84%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.474"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.77%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06194"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -12.67%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.444"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.833"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.796"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.123"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.745"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6463"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.75%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.545"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.17%  "

$ws.Range("D37").Value = "1.220.18"
$ws.Range("E37").Value = "  -1.46%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.738"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.69%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01738"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.38%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.548"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.41%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8985"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.25%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.005"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.45%  "

$ws.Range("D43").Value = "1.993.01"
$ws.Range("E43").Value = "  -0.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.49%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000115"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.98%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.504"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.64%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4567"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.34%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.576"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.84%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05512"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.55%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.428"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.74%  "
